$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# 1) Table on slide 16 (the "PLENARY - COMPLETE THE MISSING GAPS" slide):
#    switch its table style away from the custom "Table_0" style to the
#    built-in table style {E7E03FF5-F78F-4E26-A989-0BF2C1EED4EB}.
# -----------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$tableShape = $null
$shapeIndex = 1
while ($shapeIndex -le $slide16.Shapes.Count) {
    $candidate = $slide16.Shapes.Item($shapeIndex)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
    $shapeIndex = $shapeIndex + 1
}
$table = $tableShape.Table
$table.ApplyStyle("{E7E03FF5-F78F-4E26-A989-0BF2C1EED4EB}", $true)

# -----------------------------------------------------------------------
# 2) Re-colour the presentation's theme (ppt/theme/theme1.xml, used by the
#    slide master / every slide) from the custom "Integral" palette to the
#    stock Office theme palette. The font scheme and format scheme are
#    already identical between the two themes, so only the twelve scheme
#    colors need to change.
# -----------------------------------------------------------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

$dk1 = $themeColors.Colors(1)
$dk1.RGB = 0x000000

$lt1 = $themeColors.Colors(2)
$lt1.RGB = 0xFFFFFF

$dk2 = $themeColors.Colors(3)
$dk2.RGB = 0x6A5444

$lt2 = $themeColors.Colors(4)
$lt2.RGB = 0xE6E6E7

$accent1 = $themeColors.Colors(5)
$accent1.RGB = 0xD59B5B

$accent2 = $themeColors.Colors(6)
$accent2.RGB = 0x317DED

$accent3 = $themeColors.Colors(7)
$accent3.RGB = 0xA5A5A5

$accent4 = $themeColors.Colors(8)
$accent4.RGB = 0x00C0FF

$accent5 = $themeColors.Colors(9)
$accent5.RGB = 0xC47244

$accent6 = $themeColors.Colors(10)
$accent6.RGB = 0x47AD70

$hlink = $themeColors.Colors(11)
$hlink.RGB = 0xC16305

$folHlink = $themeColors.Colors(12)
$folHlink.RGB = 0x724F95
